$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E3 value to the new description text (adds a new shared string)
$ws.Range("E3").Value = "Aluslevypari  M8, NL8SP_ExcelDescCol_IMP"

# Move the active selection to E3 (was E5)
$ws.Range("E3").Select()
